$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price-tracker refresh: Price (D) and Volume(1h) (E) columns
# are plain text cells in this sheet, so we keep writing them as text.
# For Price cells whose new reading happens to look like a pure number,
# pre-set the cell to Text format so Excel does not silently convert the
# value to a numeric type on assignment (matching the original inlineStr cells).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.439.75"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "2.589.21"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "571.72"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "155.22"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -5.08%  "
$ws.Range("E9").Value = "  -7.08%  "
$ws.Range("D10").Value = "5.85"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "0.381"
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "28.18"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "3.052.79"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("E15").Value = "  -8.18%  "
$ws.Range("D16").Value = "63.190.38"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").Value = "2.585.23"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").Value = "11.97"
$ws.Range("E18").Value = "  -4.99%  "
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("D21").Value = "342.51"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D23").Value = "67.16"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").Value = "1.83"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").Value = "9.13"
$ws.Range("E26").Value = "  -5.90%  "
$ws.Range("D27").Value = "577.50"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("D33").Value = "1.72"
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "6.56"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "0.402"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "19.71"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("D39").Value = "154.41"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "41.26"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("E43").Value = "  +6.71%  "
$ws.Range("D44").Value = "155.79"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "3.91"
$ws.Range("E45").Value = "  -4.67%  "
$ws.Range("D46").Value = "23.05"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "0.0587"
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("D48").Value = "0.624"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").Value = "0.0246"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("D51").Value = "18.74"
$ws.Range("E51").Value = "  -5.35%  "
